$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.311.40'
$ws.Range('D3').Value = '1.860.68'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '242.29'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('D6').Value = '0.7001'
$ws.Range('E6').Value = '  -2.65%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '0.07875'
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('D9').Value = '0.3125'
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('D10').Value = '24.10'
$ws.Range('E10').Value = '  -3.91%  '
$ws.Range('D11').Value = '0.07807'
$ws.Range('E11').Value = '  -4.19%  '
$ws.Range('D12').Value = '1.820.51'
$ws.Range('E12').Value = '  -3.75%  '
$ws.Range('D13').Value = '5.142'
$ws.Range('E13').Value = '  -2.24%  '
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('D15').Value = '0.6938'
$ws.Range('E15').Value = '  -2.57%  '
$ws.Range('D16').Value = '6.486'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').Value = '0.000008513'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('D18').Value = '29.308.70'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('D19').Value = '248.50'
$ws.Range('E19').Value = '  -2.50%  '
$ws.Range('D20').Value = '2.115.51'
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('D21').Value = '12.96'
$ws.Range('E21').Value = '  -3.04%  '
$ws.Range('D22').Value = '0.9995'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').Value = '7.557'
$ws.Range('E23').Value = '  -3.15%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('E25').Value = '  -3.54%  '
$ws.Range('D26').Value = '160.82'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('D27').Value = '8.930'
$ws.Range('E27').Value = '  -1.81%  '
$ws.Range('D28').Value = '18.67'
$ws.Range('E28').Value = '  -2.28%  '
$ws.Range('D29').Value = '1.575'
$ws.Range('E29').Value = '  +4.29%  '
$ws.Range('D30').Value = '4.279'
$ws.Range('E30').Value = '  -3.39%  '
$ws.Range('D31').Value = '4.254'
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('D32').Value = '1.205'
$ws.Range('E32').Value = '  -1.52%  '
$ws.Range('D33').Value = '0.05236'
$ws.Range('D34').Value = '1.880'
$ws.Range('E34').Value = '  -3.83%  '
$ws.Range('D35').Value = '0.7507'
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('D37').Value = '2.695'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('E38').Value = '  -1.88%  '
$ws.Range('D39').Value = '1.262.00'
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('D40').Value = '2.744'
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('D41').Value = '0.8998'
$ws.Range('E41').Value = '  -0.72%  '
$ws.Range('D42').Value = '111.07'
$ws.Range('E42').Value = '  -1.84%  '
$ws.Range('D43').Value = '5.950'
$ws.Range('E43').Value = '  -8.17%  '
$ws.Range('D44').Value = '69.42'
$ws.Range('E44').Value = '  -6.88%  '
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '2.013.20'
$ws.Range('E46').Value = '  -1.20%  '
$ws.Range('E47').Value = '  -4.32%  '
$ws.Range('D48').Value = '9.548'
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('D49').Value = '0.5180'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('D50').Value = '1.778'
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('D51').Value = '0.4263'
$ws.Range('E51').Value = '  -2.62%  '
